# Fonds de solidarite - add 2020-09-01 data
# This sheet stores every value as text (inlineStr), including numbers like
# "190" or "448016.00" and zero-padded codes like "00"/"06". To keep them as
# text (not get auto-converted to numbers, which would strip trailing/leading
# zeros), every write below first forces the target cell(s) to text format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- 1) Update nombre_aides (C) / montant_total (D) for rows that are not
#        affected by the later row insertion (rows 2-75, i.e. before Mayotte). ---
$updates1 = @(
    @{ Row = 2;  C = "192";  D = "452016.00" },
    @{ Row = 3;  C = "1028"; D = "3295120.33" },
    @{ Row = 8;  C = "43";   D = "86000.00" },
    @{ Row = 16; C = "429";  D = "1366161.23" },
    @{ Row = 22; C = "328";  D = "971679.20" },
    @{ Row = 33; C = "111";  D = "321173.00" },
    @{ Row = 34; C = "587";  D = "1936276.10" },
    @{ Row = 35; C = "246";  D = "1293520.11" },
    @{ Row = 38; C = "24";   D = "52200.00" },
    @{ Row = 39; C = "40";   D = "105536.00" },
    @{ Row = 40; C = "173";  D = "470779.00" },
    @{ Row = 41; C = "84";   D = "294298.00" },
    @{ Row = 44; C = "65";   D = "167183.00" },
    @{ Row = 45; C = "29";   D = "106621.84" },
    @{ Row = 46; C = "93";   D = "413774.61" },
    @{ Row = 47; C = "52";   D = "303203.00" },
    @{ Row = 50; C = "16";   D = "35850.00" },
    @{ Row = 75; C = "13";   D = "68000.00" }
)

foreach ($u in $updates1) {
    Set-TextValue $ws.Range("C" + $u.Row) $u.C
    Set-TextValue $ws.Range("D" + $u.Row) $u.D
}

# --- 2) Insert a brand-new row at 76 (Mayotte / "1 ou 2 salariés"), pushing
#        the former rows 76-105 down to 77-106. ---
$ws.Rows.Item(76).Insert()

Set-TextValue $ws.Range("A76") "Fonds de solidarité"
Set-TextValue $ws.Range("B76") "VOLET2"
Set-TextValue $ws.Range("C76") "3"
Set-TextValue $ws.Range("D76") "14000.00"
Set-TextValue $ws.Range("E76") "06"
Set-TextValue $ws.Range("F76") "Mayotte"
Set-TextValue $ws.Range("G76") "01"
Set-TextValue $ws.Range("H76") "1 ou 2 salariés"

# --- 3) Update nombre_aides (C) / montant_total (D) for the rows that moved
#        during the insertion (former rows 82/83/85/86, now 83/84/86/87). ---
$updates2 = @(
    @{ Row = 83; C = "232"; D = "597326.09" },
    @{ Row = 84; C = "903"; D = "2909509.26" },
    @{ Row = 86; C = "119"; D = "591984.52" },
    @{ Row = 87; C = "30";  D = "184080.04" }
)

foreach ($u in $updates2) {
    Set-TextValue $ws.Range("C" + $u.Row) $u.C
    Set-TextValue $ws.Range("D" + $u.Row) $u.D
}
